$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was 45203 (2023-10-04)
# for every data row (2-498) and must become 45204 (2023-10-05).
$range = $ws.Range("C2:C498")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45204
    }
}
